# Applies the update described by the commit:
# "added functionality to provide a project name when creating a project"
#
# Effect on the data:
#   - All Jira issue keys in column A (rows 2-39) are renumbered, shifting
#     the numeric suffix up by 106 (e.g. TD-7102 -> TD-7208).
#   - Every cell (in columns B or D) whose value is exactly "Pilot" (the
#     project name) is renamed to "PILOTWIL".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)  # column A - IssueKey
    $val = $cell.Value2
    if ($val -match '^(TD-)(\d+)$') {
        $prefix = $Matches[1]
        $num = [int]$Matches[2]
        $newNum = $num + 106
        $cell.Value = "$prefix$newNum"
    }
}

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in 2, 4) {  # column B (Summary) and column D (ParentSummary)
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.Value2 -eq 'Pilot') {
            $cell.Value = 'PILOTWIL'
        }
    }
}
